$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.026.25'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.897.59'
$ws.Range("E3").Value = '  -1.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7367'
$ws.Range("E5").Value = '  -2.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.88'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9994'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3096'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.40'
$ws.Range("E9").Value = '  -4.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06900'
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7714'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07957'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '1.897.10'
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.230'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.59'
$ws.Range("E15").Value = '  -2.84%  '
$ws.Range("D16").Value = '30.029.23'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.15'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007781'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.71'
$ws.Range("E20").Value = '  -4.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9989'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '2.152.19'
$ws.Range("E22").Value = '  -1.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.894'
$ws.Range("E24").Value = '  +3.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.320'
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.38'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.82'
$ws.Range("E27").Value = '  -0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1276'
$ws.Range("E28").Value = '  -4.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.040'
$ws.Range("E29").Value = '  -7.50%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.539'
$ws.Range("E31").Value = '  +1.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.300'
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05112'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.284'
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7371'
$ws.Range("E36").Value = '  -1.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  -2.09%  '
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.798'
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.308'
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.35'
$ws.Range("E41").Value = '  -4.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4459'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.940'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.0000'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8365'
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.740'
$ws.Range("E46").Value = '  +4.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.21'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.850'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").Value = '2.053.96'
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.67'
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '933.62'
$ws.Range("E51").Value = '  -5.22%  '
